$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-26: 2023-09-15 -> 2023-09-16
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
